$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a new row before row 30; this shifts the existing row 30
# (and everything below it) down to row 31, carrying along its
# values and formatting.
$ws.Rows.Item(30).Insert()

# Populate the newly inserted row 30 with the new weekly record.
$ws.Range("A30").Value = 2
$ws.Range("B30").Value = "Comercializadora del Agro de Limarí"
$ws.Range("C30").Value = "Coquimbo"
$ws.Range("D30").Value = 44448
$ws.Range("D30").NumberFormat = $ws.Range("D31").NumberFormat
$ws.Range("E30").Value = 4
$ws.Range("F30").Value = 100112026
$ws.Range("G30").Value = "Haba"
$ws.Range("H30").Value = "Sin especificar"
$ws.Range("I30").Value = "Primera"
$ws.Range("J30").Value = 800
$ws.Range("K30").Value = 10000
$ws.Range("L30").Value = 12000
$ws.Range("M30").Value = 11000
$ws.Range("N30").Value = "$/saco 25 kilos"
$ws.Range("O30").Value = "Provincia de Limarí"
$ws.Range("P30").Value = 440
$ws.Range("Q30").Value = 25
$ws.Range("R30").Value = "Hortaliza"

$wb.Save()
